$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the current row 20 ("Haba"
# data for Vega Modelo de Temuco), pushing the former rows 20-27 down to
# 21-28. Insert a fresh row at position 20 so the rest of the sheet shifts
# down, matching the rest of the table's layout/style.
$ws.Rows("20:20").Insert()

# Fill in the newly inserted row 20 with the new weekly record.
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44468
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100112026
$ws.Cells.Item(20, 7).Value = "Haba"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 14000
$ws.Cells.Item(20, 13).Value = 14000
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 560
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
